{"js": "// Move the \"_GoBack\" bookmark from the start of the \"\u0418\u043c\u044f \u2013 \u043c\u0435\u043d\u0435\u0434\u0436\u0435\u0440_\u0438\u043c\u044f\"\n// paragraph down to the start of the \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0434\u0430\u0442\u0430\"\n// paragraph, and right-align the two report-footer paragraphs\n// (\"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430 ...\" and \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430 ...\").\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1) Remove the old \"_GoBack\" bookmark (currently at the very start of\n//    the \"\u0418\u043c\u044f \u2013 \u043c\u0435\u043d\u0435\u0434\u0436\u0435\u0440_\u0438\u043c\u044f\" paragraph).\ndoc.deleteBookmark(\"_GoBack\");\n\n// 2) Right-align the \"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" paragraph.\nconst timeResults = body.search(\"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\", { matchCase: false });\ntimeResults.load(\"items\");\nawait context.sync();\nif (timeResults.items.length > 0) {\n  const timeParagraph = timeResults.items[0].paragraphs.getFirst();\n  timeParagraph.alignment = \"Right\";\n}\n\n// 3) Right-align the \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" paragraph and re-insert\n//    the \"_GoBack\" bookmark collapsed at its very start (before the\n//    \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" run).\nconst dateResults = body.search(\"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\", { matchCase: false });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  const dateParagraph = dateResults.items[0].paragraphs.getFirst();\n  dateParagraph.alignment = \"Right\";\n  const dateStart = dateParagraph.getRange(\"Start\");\n  dateStart.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Word COM script: move the \"_GoBack\" bookmark from the start of the\n# \"\u0418\u043c\u044f \u2013 \u043c\u0435\u043d\u0435\u0434\u0436\u0435\u0440_\u0438\u043c\u044f\" paragraph down to the start of the\n# \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430: \u0434\u0430\u0442\u0430\" paragraph, and right-align the two\n# report-footer paragraphs (\"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430 ...\" and\n# \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430 ...\").\n\n$d = $word.ActiveDocument\n\n# 1) Remove the old \"_GoBack\" bookmark (it currently sits at the very\n#    start of the \"\u0418\u043c\u044f \u2013 \u043c\u0435\u043d\u0435\u0434\u0436\u0435\u0440_\u0438\u043c\u044f\" paragraph).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Right-align the \"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" paragraph.\n$timeRng = $d.Content\n$timeFound = $timeRng.Find.Execute(\"\u0412\u0440\u0435\u043c\u044f \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\")\nif ($timeFound) {\n    $timeRng.Paragraphs(1).Alignment = 2\n}\n\n# 3) Right-align the \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" paragraph and add the\n#    \"_GoBack\" bookmark back, collapsed at the very start of that\n#    paragraph (before the \"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\" run).\n$dateRng = $d.Content\n$dateFound = $dateRng.Find.Execute(\"\u0414\u0430\u0442\u0430 \u0444\u043e\u0440\u043c\u0438\u0440\u043e\u0432\u0430\u043d\u0438\u044f \u043e\u0442\u0447\u0435\u0442\u0430\")\nif ($dateFound) {\n    $dateParagraph = $dateRng.Paragraphs(1)\n    $dateParagraph.Alignment = 2\n    $startPoint = $d.Range($dateParagraph.Range.Start, $dateParagraph.Range.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $startPoint)\n}\n"}
